$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "254.77"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.68%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.91"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-7.31%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.222"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.27%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05872"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.94%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.716"
$ws.Range("D6").ClearFormats()
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8634"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.63%"
$ws.Range("E7").ClearFormats()
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9680"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "12.78%"
$ws.Range("E8").ClearFormats()
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.01055"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.90%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1410"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.95%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07160"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.09%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03181"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.72%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09221"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.49%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001541"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.80%"
$ws.Range("E14").ClearFormats()
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005802"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.61%"
$ws.Range("E15").ClearFormats()
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.501"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.29%"
$ws.Range("E16").ClearFormats()
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.218"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.97%"
$ws.Range("E17").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.21%"
$ws.Range("E18").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.21%"
$ws.Range("E19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03471"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.44%"
$ws.Range("E20").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.556"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.11%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04166"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.36%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1380"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.13%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001227"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.13%"
$ws.Range("E25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004796"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "15.30%"
$ws.Range("E26").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.02%"
$ws.Range("E27").ClearFormats()
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "1.16%"
$ws.Range("E28").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03811"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.48%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005749"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.21%"
$ws.Range("E41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1099"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.84%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002342"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.48%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01075"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "13.75%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005243"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.57%"
$ws.Range("E45").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.02%"
$ws.Range("E46").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "40.82%"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002132"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-13.58%"
$ws.Range("E48").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.02%"
$ws.Range("E49").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.02%"
$ws.Range("E50").ClearFormats()
